# Update the "Förändrad" (Changed) date column (C) for rows 2-12
# from 2023-10-05 (45204) to 2023-10-08 (45207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
